# This script re-shuffles the per-row stimulus data (columns H, I, K, L, M, N,
# O, P, Q, R, S, T, U, V) across rows 2-41 of the active sheet according to a
# fixed row permutation. Rows 3 and 14 keep their original content. Columns
# A-G and J are left untouched (they already contain the sequential trial
# bookkeeping data that does not change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) whose content is permuted between rows.
$cols = @("H", "I", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V")

$firstRow = 2
$lastRow = 41

# Snapshot the "before" values for every affected column/row so that later
# writes don't clobber data we still need to read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row (the content that ends up in the destination
# row is exactly what used to be in the source row). Rows not listed map to
# themselves (i.e. unchanged).
$mapping = @{
    2 = 29; 4 = 36; 5 = 23; 6 = 27; 7 = 24; 8 = 30; 9 = 28; 10 = 15;
    11 = 20; 12 = 4; 13 = 6; 15 = 13; 16 = 34; 17 = 40; 18 = 19; 19 = 37;
    20 = 31; 21 = 25; 22 = 41; 23 = 33; 24 = 17; 25 = 7; 26 = 10; 27 = 12;
    28 = 38; 29 = 21; 30 = 18; 31 = 22; 32 = 5; 33 = 9; 34 = 8; 35 = 11;
    36 = 39; 37 = 2; 38 = 16; 39 = 32; 40 = 35; 41 = 26
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $r
    if ($mapping.ContainsKey($r)) {
        $srcRow = $mapping[$r]
    }
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $srcVals[$col]
    }
}
